# Apply the edits described by the commit:
#   "Finish with GP16 state chart / WIP: Need to do traceability chart"
#
# 1) The cached "Header & Footer" fixed-date text was refreshed from
#    2021-03-27 -> 2021-04-04 everywhere it is cached (slide master,
#    every slide layout, and the notes master).
# 2) On slide 2 ("PR-7c ...") the alarm-time-advance text was reworded
#    from "... at 5 minutes of alarm time ..." to
#    "... at 1 hour of alarm time ...".

$p = $ppt.ActivePresentation

$oldDate = "2021-03-27"
$newDate = "2021-04-04"

function Update-DateIfMatch {
    param($shape)

    if (-not $shape.HasTextFrame) { return }
    if (-not $shape.TextFrame.HasText) { return }

    $tr = $shape.TextFrame.TextRange
    if ($tr.Text -eq $oldDate) {
        $tr.Text = $newDate
    }
}

# --- Slide master date placeholder ---
for ($j = 1; $j -le $p.SlideMaster.Shapes.Count; $j++) {
    Update-DateIfMatch($p.SlideMaster.Shapes.Item($j))
}

# --- Every slide layout's date placeholder ---
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        Update-DateIfMatch($layout.Shapes.Item($j))
    }
}

# --- Notes master date placeholder ---
$notesMaster = $p.NotesMaster
for ($j = 1; $j -le $notesMaster.Shapes.Count; $j++) {
    Update-DateIfMatch($notesMaster.Shapes.Item($j))
}

# --- Slide 2: reword the PR-7c bullet ---
$slide2 = $p.Slides.Item(2)
$contentShape = $slide2.Shapes.Item(2)
$bodyRange = $contentShape.TextFrame.TextRange

# The PR-7c bullet is the 7th paragraph in the content placeholder.
$pr7c = $bodyRange.Paragraphs(7, 1)

$target = "5 minutes"
$full = $pr7c.Text
$idx = $full.IndexOf($target)
if ($idx -ge 0) {
    $sub = $pr7c.Characters($idx + 1, $target.Length)
    $sub.Text = "1 hour"
}
